$d = $word.ActiveDocument

# Update the trailing whitespace after "Senior Software Engineer" so the
# job title run shrinks from 66 to 53 trailing spaces (13 fewer spaces).
$oldTitle = "Senior Software Engineer                                                                  "
$newTitle = "Senior Software Engineer                                                     "
$d.Content.Find.Execute($oldTitle, $false, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)

# Update Gojek exit date from "Present" to "November 2022".
$oldDate = [char]0x2013 + " Present"
$newDate = [char]0x2013 + " November 2022"
$d.Content.Find.Execute($oldDate, $false, $false, $false, $false, $false, $true, 1, $false, $newDate, 2)
